# [Kadastro App] Yeni kayit eklendi: 2964
# Append a new record row (row 32) to the "Kayitlar" (overview) sheet and
# to the "Erdemli" district sheet, which mirrors the same table.

$wb = $excel.ActiveWorkbook

$newRowNumber = 32
$kayitNo      = "2964"
$tarih        = "2025-09-10"
$birim        = "Erdemli"
$parselSayisi = "1"
$is           = "ÇAP"
$personeller  = "AYHAN KARADAYI (K.Teknisyeni)"

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $targetRow = $ws.Range("A" + $newRowNumber + ":F" + $newRowNumber)

    # Leading apostrophes force these numeric/date-looking values to be
    # stored as literal text (matching the rest of the table, where every
    # column - including "numeric" ones like Kayit No / Parsel Sayisi - is
    # kept as text and Excel's "number stored as text" warning is
    # suppressed via ignoredErrors). ClearFormats() afterwards drops the
    # implicit quote-prefix style so the new cells keep the sheet's default
    # (unstyled) formatting, just like every other cell in the table.
    $ws.Range("A" + $newRowNumber).Value = "'" + $kayitNo
    $ws.Range("B" + $newRowNumber).Value = "'" + $tarih
    $ws.Range("C" + $newRowNumber).Value = $birim
    $ws.Range("D" + $newRowNumber).Value = "'" + $parselSayisi
    $ws.Range("E" + $newRowNumber).Value = $is
    $ws.Range("F" + $newRowNumber).Value = $personeller

    $targetRow.ClearFormats()
}
